# Merged CDS test suites to create CDS_Regression suite.
# This updates the "Files" query text so the Participant ID / Sample ID
# coalesce() calls fall back to an empty string instead of the literal
# "Not specified in data" placeholder, and resets the sheet's selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# Update the query text stored in B4 (Files tab query).
[string]$filesQuery = $ws.Range("B4").Value2
$filesQuery = $filesQuery.Replace(
    "coalesce(p.participant_id, 'Not specified in data')",
    "coalesce(p.participant_id, '')"
)
$filesQuery = $filesQuery.Replace(
    "coalesce(samp.sample_id, 'Not specified in data')",
    "coalesce(samp.sample_id, '')"
)
$ws.Range("B4").Value2 = $filesQuery

# Reset the active sheet view/selection to B2 (also clears the A3 scroll anchor).
$ws.Activate()
$ws.Range("B2").Select()
